# Update the cached "datetimeFigureOut" date field text from 3/23/2022 to
# 4/26/2022 everywhere it appears: on the Slide Master and on every Slide
# Layout's "Date Placeholder" shape.

$p = $ppt.ActivePresentation
$oldDate = "3/23/2022"
$newDate = "4/26/2022"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq 16) {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout (CustomLayouts) under the (only) Design
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

Write-Output "done"
